# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the 632c8bd2-d307-4447-979d-dc53acf932fd file across the Overview, zh-cn and
# de-de sheets, reflecting a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 4 = 632c8bd2... file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-11-29 03:59:49"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 4 = 632c8bd2... file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-11-29 03:59:36"

# de-de sheet: column H = "Latest Handoff Datetime", row 4 = 632c8bd2... file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-11-29 03:59:49"
